$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.573375
$ws.Range("H2").Value = 28.720125
$ws.Range("I2").Value = 0.1037691388643484
$ws.Range("J2").Value = 0.1037691388643484
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 1394.84786647325
$ws.Range("R2").Value = 12553.63079825925
$ws.Range("S2").Value = 0.02973958923673114
$ws.Range("T2").Value = 0.02973958923673114
$ws.Range("G3").Value = 9.573375
$ws.Range("H3").Value = 28.720125
$ws.Range("I3").Value = 0.1037691388643484
$ws.Range("J3").Value = 0.1037691388643484
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 1615.9828535165
$ws.Range("R3").Value = 14543.8456816485
$ws.Range("S3").Value = 0.03445441430017274
$ws.Range("T3").Value = 0.03445441430017274
$ws.Range("G4").Value = 9.573375
$ws.Range("H4").Value = 28.720125
$ws.Range("I4").Value = 0.1037691388643484
$ws.Range("J4").Value = 0.1037691388643484
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 1226.599725932
$ws.Range("R4").Value = 11039.397533388
$ws.Range("S4").Value = 0.02615236606364645
$ws.Range("T4").Value = 0.02615236606364645
$ws.Range("G5").Value = 9.573375
$ws.Range("H5").Value = 28.720125
$ws.Range("I5").Value = 0.1037691388643484
$ws.Range("J5").Value = 0.1037691388643484
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 629.555469671625
$ws.Range("R5").Value = 5665.999227044625
$ws.Range("S5").Value = 0.01342276926379808
$ws.Range("T5").Value = 0.01342276926379808
$ws.Range("G6").Value = 47.94465366666667
$ws.Range("I6").Value = 0.5196887643218222
$ws.Range("J6").Value = 0.5196887643218222
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 6985.571742018766
$ws.Range("R6").Value = 62870.14567816889
$ws.Range("S6").Value = 0.1489395647975769
$ws.Range("T6").Value = 0.1489395647975769
$ws.Range("G7").Value = 47.94465366666667
$ws.Range("I7").Value = 0.5196887643218222
$ws.Range("J7").Value = 0.5196887643218222
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("Q7").Value = 8093.043283389646
$ws.Range("R7").Value = 72837.38955050681
$ws.Range("S7").Value = 0.1725519956033927
$ws.Range("T7").Value = 0.1725519956033927
$ws.Range("G8").Value = 47.94465366666667
$ws.Range("I8").Value = 0.5196887643218222
$ws.Range("J8").Value = 0.5196887643218222
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 6142.964111135102
$ws.Range("R8").Value = 55286.67700021591
$ws.Range("S8").Value = 0.1309743046193652
$ws.Range("T8").Value = 0.1309743046193652
$ws.Range("G9").Value = 47.94465366666667
$ws.Range("I9").Value = 0.5196887643218222
$ws.Range("J9").Value = 0.5196887643218222
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 3152.892157401307
$ws.Range("R9").Value = 28376.02941661176
$ws.Range("S9").Value = 0.06722289930148746
$ws.Range("T9").Value = 0.06722289930148746
$ws.Range("G10").Value = 11.32006633333333
$ws.Range("H10").Value = 33.960199
$ws.Range("I10").Value = 0.122702133291269
$ws.Range("J10").Value = 0.122702133291269
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 1649.34209444273
$ws.Range("R10").Value = 14844.07884998457
$ws.Range("S10").Value = 0.03516566758179666
$ws.Range("T10").Value = 0.03516566758179666
$ws.Range("G11").Value = 11.32006633333333
$ws.Range("H11").Value = 33.960199
$ws.Range("I11").Value = 0.122702133291269
$ws.Range("J11").Value = 0.122702133291269
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 1910.823831233611
$ws.Range("R11").Value = 17197.4144811025
$ws.Range("S11").Value = 0.04074072679218186
$ws.Range("T11").Value = 0.04074072679218186
$ws.Range("G12").Value = 11.32006633333333
$ws.Range("H12").Value = 33.960199
$ws.Range("I12").Value = 0.122702133291269
$ws.Range("J12").Value = 0.122702133291269
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 1450.396569861593
$ws.Range("R12").Value = 13053.56912875434
$ws.Range("S12").Value = 0.03092394465004174
$ws.Range("T12").Value = 0.03092394465004174
$ws.Range("G13").Value = 11.32006633333333
$ws.Range("H13").Value = 33.960199
$ws.Range("I13").Value = 0.122702133291269
$ws.Range("J13").Value = 0.122702133291269
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 744.4197764315737
$ws.Range("R13").Value = 6699.777987884163
$ws.Range("S13").Value = 0.01587179426724871
$ws.Range("T13").Value = 0.01587179426724871
$ws.Range("G14").Value = 23.41838033333333
$ws.Range("H14").Value = 70.25514099999999
$ws.Range("I14").Value = 0.2538399635225604
$ws.Range("J14").Value = 0.2538399635225604
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 3412.075453453888
$ws.Range("R14").Value = 30708.679081085
$ws.Range("S14").Value = 0.07274895339447961
$ws.Range("T14").Value = 0.07274895339447961
$ws.Range("G15").Value = 23.41838033333333
$ws.Range("H15").Value = 70.25514099999999
$ws.Range("I15").Value = 0.2538399635225604
$ws.Range("J15").Value = 0.2538399635225604
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 3953.015637201581
$ws.Range("R15").Value = 35577.14073481423
$ws.Range("S15").Value = 0.08428235374083685
$ws.Range("T15").Value = 0.08428235374083685
$ws.Range("G16").Value = 23.41838033333333
$ws.Range("H16").Value = 70.25514099999999
$ws.Range("I16").Value = 0.2538399635225604
$ws.Range("J16").Value = 0.2538399635225604
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 3000.507020631492
$ws.Range("R16").Value = 27004.56318568342
$ws.Range("S16").Value = 0.06397389166255704
$ws.Range("T16").Value = 0.06397389166255704
$ws.Range("G17").Value = 23.41838033333333
$ws.Range("H17").Value = 70.25514099999999
$ws.Range("I17").Value = 0.2538399635225604
$ws.Range("J17").Value = 0.2538399635225604
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 1540.017959152379
$ws.Range("R17").Value = 13860.16163237142
$ws.Range("S17").Value = 0.03283476472468697
$ws.Range("T17").Value = 0.03283476472468697
